# "spell check in presentation"
#
# Slide 2 ("Executive Summary") has a single body placeholder (Shapes.Item(1))
# containing three paragraphs. The middle paragraph is re-worded slightly
# ("for both event attendees and organizers alike" -> "for event attendees
# and organizers"), and PowerPoint marks all three runs/paragraph ends as
# spell-checked ("dirty") after the edit.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# --- Paragraph 1 (unchanged text, just re-affirm / mark as clean) ---------
$run1 = $tr.Paragraphs(1, 1).Runs(1, 1)
$run1.Text = "TKT is a ticket sale and management application built with ERC721 NFTs on the Ethereum Network."
$run1.Font.Dirty = 0

# --- Paragraph 2 (actual wording fix) --------------------------------------
$run2 = $tr.Paragraphs(2, 1).Runs(1, 1)
$run2.Text = "TKT offers a secure and transparent platform for event attendees and organizers. Because TKTs are minted on the Ethereum blockchain, transactions are immutable and attendees have total ownership of TKTs in their wallet, while organizers can guarantee authenticity."
$run2.Font.Dirty = 0

# --- Paragraph 3 (unchanged text, just re-affirm / mark as clean) ---------
$run3 = $tr.Paragraphs(3, 1).Runs(1, 1)
$run3.Text = "TKT aims to decentralize and disrupt the `$60+ billion global ticketing industry."
$run3.Font.Dirty = 0

# --- Notes page for this slide: the slide-image placeholder shifts by
# 300 EMU (0.03 pt) on the X axis, an artifact PowerPoint leaves behind
# when it re-lays-out the notes page after a save.
$notesShape = $s.NotesPage.Shapes.Item(1)
$notesShape.Left = 30.0
